$wb = $excel.ActiveWorkbook

# This script applies a scheduled-runner data refresh to the per-job (Leve)
# profit tables on each crafting-class sheet (ALC, ARM, BSM, CRP, CUL, GSM,
# LTW, WVR). Columns H:N are cached marketboard-derived numbers:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
# Each block below updates one row to the freshly-fetched values.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 87: There Was a Late Fee / Noble Gold
$ws.Range("H87").Value = 40460.91
$ws.Range("J87").Value = 40460.91
$ws.Range("L87").Value = 40460.91
$ws.Range("N87").Value = -42956.91

# row 90: A Gate Arcane Is Dragon's Bane (L) / Noble Gold
$ws.Range("H90").Value = 40460.91
$ws.Range("J90").Value = 40460.91
$ws.Range("L90").Value = 121382.73
$ws.Range("N90").Value = -133862.73

# row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 758.0476
$ws.Range("I107").Value = 714.9375
$ws.Range("J107").Value = 896
$ws.Range("K107").Value = 714.9375
$ws.Range("L107").Value = 896
$ws.Range("M107").Value = 1205.0625
$ws.Range("N107").Value = -4736

# row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1319.5873
$ws.Range("I132").Value = 1018.9
$ws.Range("J132").Value = 7333.3335
$ws.Range("K132").Value = 3056.7
$ws.Range("L132").Value = 22000.0005
$ws.Range("M132").Value = -526.6999999999998
$ws.Range("N132").Value = -27060.0005

# row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1037483.1
$ws.Range("I138").Value = 2501.0476
$ws.Range("J138").Value = 1361880.5
$ws.Range("K138").Value = 7503.1428
$ws.Range("L138").Value = 4085641.5
$ws.Range("M138").Value = -2363.1428
$ws.Range("N138").Value = -4095921.5


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 18: Still the Best / Brass Alembic
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 12692.4
$ws.Range("I32").Value = 10298.526
$ws.Range("J32").Value = 29748.75
$ws.Range("K32").Value = 10298.526
$ws.Range("L32").Value = 29748.75
$ws.Range("M32").Value = -10011.526
$ws.Range("N32").Value = -30322.75

# row 140: A Hand for a Deckhand / Ra'Kaznar Gloves of Scouting
$ws.Range("H140").Value = 39692.4
$ws.Range("J140").Value = 41123.54
$ws.Range("L140").Value = 41123.54
$ws.Range("N140").Value = -51483.54


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 80: Unbreaker / Titanium Ingot
$ws.Range("H80").Value = 244.1923
$ws.Range("I80").Value = 130.85715
$ws.Range("J80").Value = 285.94736
$ws.Range("K80").Value = 130.85715
$ws.Range("L80").Value = 285.94736
$ws.Range("M80").Value = 867.14285
$ws.Range("N80").Value = -2281.94736

# row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Range("H83").Value = 244.1923
$ws.Range("I83").Value = 130.85715
$ws.Range("J83").Value = 285.94736
$ws.Range("K83").Value = 654.28575
$ws.Range("L83").Value = 1429.7368
$ws.Range("M83").Value = 4337.71425
$ws.Range("N83").Value = -11413.7368


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1888.75
$ws.Range("I16").Value = 1850
$ws.Range("J16").Value = 2005
$ws.Range("K16").Value = 1850
$ws.Range("L16").Value = 2005
$ws.Range("M16").Value = -1563
$ws.Range("N16").Value = -2579

# row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 9070.571
$ws.Range("J31").Value = 9070.571
$ws.Range("L31").Value = 9070.571
$ws.Range("N31").Value = -9660.571

# row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 9070.571
$ws.Range("J34").Value = 9070.571
$ws.Range("L34").Value = 9070.571
$ws.Range("N34").Value = -9474.571

# row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 2302.2222
$ws.Range("I99").Value = 2436.6667
$ws.Range("J99").Value = 2033.3334
$ws.Range("K99").Value = 2436.6667
$ws.Range("L99").Value = 2033.3334
$ws.Range("M99").Value = -938.6667000000002
$ws.Range("N99").Value = -5029.3334

# row 109: Playing the Market / White Oak Necklace
$ws.Range("H109").Value = 29000
$ws.Range("I109").Value = 25000
$ws.Range("J109").Value = 33000
$ws.Range("K109").Value = 25000
$ws.Range("L109").Value = 33000
$ws.Range("M109").Value = -23960
$ws.Range("N109").Value = -35080

# row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1888.75
$ws.Range("I113").Value = 1850
$ws.Range("J113").Value = 2005
$ws.Range("K113").Value = 1850
$ws.Range("L113").Value = 2005
$ws.Range("M113").Value = 320
$ws.Range("N113").Value = -6345

# row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 2302.2222
$ws.Range("I126").Value = 2436.6667
$ws.Range("J126").Value = 2033.3334
$ws.Range("K126").Value = 7310.000100000001
$ws.Range("L126").Value = 6100.0002
$ws.Range("M126").Value = -4840.000100000001
$ws.Range("N126").Value = -11040.0002


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 754.7959
$ws.Range("I113").Value = 770.25
$ws.Range("J113").Value = 725.7059
$ws.Range("K113").Value = 2310.75
$ws.Range("L113").Value = 2177.1177
$ws.Range("M113").Value = -140.75
$ws.Range("N113").Value = -6517.117700000001

# row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 1009.3333
$ws.Range("J117").Value = 1009.3333
$ws.Range("L117").Value = 3027.9999
$ws.Range("N117").Value = -9911.999899999999

# row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 2344.9
$ws.Range("I129").Value = 4143.3335
$ws.Range("J129").Value = 1574.1428
$ws.Range("K129").Value = 12430.0005
$ws.Range("L129").Value = 4722.428400000001
$ws.Range("M129").Value = -7430.000499999998
$ws.Range("N129").Value = -14722.4284

# row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 45460376
$ws.Range("I137").Value = 55559676
$ws.Range("K137").Value = 166679028
$ws.Range("M137").Value = -166673928


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 5942.4585
$ws.Range("I102").Value = 5638.933
$ws.Range("J102").Value = 6448.3335
$ws.Range("K102").Value = 5638.933
$ws.Range("L102").Value = 6448.3335
$ws.Range("M102").Value = -4016.933
$ws.Range("N102").Value = -9692.333500000001

# row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3284.0557
$ws.Range("I126").Value = 1961.25
$ws.Range("J126").Value = 4342.3
$ws.Range("K126").Value = 5883.75
$ws.Range("L126").Value = 13026.9
$ws.Range("M126").Value = -3413.75
$ws.Range("N126").Value = -17966.9


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 6985.3076
$ws.Range("I7").Value = 3988
$ws.Range("K7").Value = 3988
$ws.Range("M7").Value = -3876

# row 38: Emergency Patches / Skull Eyepatch
$ws.Range("H38").Value = 19999.334
$ws.Range("J38").Value = 19999.334
$ws.Range("L38").Value = 19999.334
$ws.Range("N38").Value = -20819.334

# row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4071.2942
$ws.Range("I40").Value = 4016.3076
$ws.Range("J40").Value = 4250
$ws.Range("K40").Value = 4016.3076
$ws.Range("L40").Value = 4250
$ws.Range("M40").Value = -3880.3076
$ws.Range("N40").Value = -4522

# row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6829.9
$ws.Range("I122").Value = 4614.143
$ws.Range("J122").Value = 12000
$ws.Range("K122").Value = 13842.429
$ws.Range("L122").Value = 36000
$ws.Range("M122").Value = -11392.429
$ws.Range("N122").Value = -40900

# row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 6985.3076
$ws.Range("I126").Value = 3988
$ws.Range("K126").Value = 11964
$ws.Range("M126").Value = -9494

# row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 7572.4287
$ws.Range("I132").Value = 13584.667
$ws.Range("J132").Value = 3063.25
$ws.Range("K132").Value = 40754.001
$ws.Range("L132").Value = 9189.75
$ws.Range("M132").Value = -38224.001
$ws.Range("N132").Value = -14249.75

# row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 6818.516
$ws.Range("I136").Value = 5784.923
$ws.Range("J136").Value = 7565
$ws.Range("K136").Value = 17354.769
$ws.Range("L136").Value = 22695
$ws.Range("M136").Value = -14804.769
$ws.Range("N136").Value = -27795


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 109: Turban in Training / Brightlinen Turban of Crafting
$ws.Range("H109").Value = 62974
$ws.Range("J109").Value = 62974
$ws.Range("L109").Value = 62974
$ws.Range("N109").Value = -65748

# row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 7368
$ws.Range("I122").Value = 1884.5
$ws.Range("J122").Value = 18335
$ws.Range("K122").Value = 5653.5
$ws.Range("L122").Value = 55005
$ws.Range("M122").Value = -3203.5
$ws.Range("N122").Value = -59905

# row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 44285.43
$ws.Range("I123").Value = 50000
$ws.Range("J123").Value = 43333
$ws.Range("K123").Value = 50000
$ws.Range("L123").Value = 43333
$ws.Range("M123").Value = -45100
$ws.Range("N123").Value = -53133

# row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 5226.3438
$ws.Range("I136").Value = 2989.875
$ws.Range("K136").Value = 8969.625
$ws.Range("M136").Value = -6419.625

